# Auto-generated edit script: updates crypto price/volume table
# matching the commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.545.41"
$ws.Range("E2").Value = "  -0.17%  "

$ws.Range("D3").Value = "3.233.80"
$ws.Range("E3").Value = "  +1.22%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("E5").Value = "  +0.34%  "

$ws.Range("D6").Value = "'156.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.48%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "3.232.25"
$ws.Range("E8").Value = "  +1.17%  "

$ws.Range("E9").Value = "  -0.87%  "

$ws.Range("E10").Value = "  +1.64%  "

$ws.Range("D11").Value = "'5.75"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.18%  "

$ws.Range("E12").Value = "  -2.20%  "

$ws.Range("E13").Value = "  +1.84%  "

$ws.Range("D14").Value = "'38.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.75%  "

$ws.Range("D15").Value = "3.747.33"
$ws.Range("E15").Value = "  +0.64%  "

$ws.Range("D16").Value = "66.602.89"
$ws.Range("E16").Value = "  -0.16%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'7.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.95%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.205.95"
$ws.Range("E18").Value = "  +0.25%  "

$ws.Range("E19").Value = "  +1.36%  "

$ws.Range("D20").Value = "'506.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.73%  "

$ws.Range("D21").Value = "'15.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.95%  "

$ws.Range("D22").Value = "'0.739"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("D23").Value = "'7.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.23%  "

$ws.Range("D24").Value = "'14.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.86%  "

$ws.Range("D25").Value = "'86.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.20%  "

$ws.Range("E26").Value = "  +85.77%  "

$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").Value = "'2.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.54%  "

$ws.Range("D29").Value = "'9.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.89%  "

$ws.Range("D30").Value = "'2.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.85%  "

$ws.Range("D31").Value = "'2.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.27%  "

$ws.Range("D32").Value = "'6.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.19%  "

$ws.Range("D33").Value = "'28.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("E35").Value = "  -4.99%  "

$ws.Range("D36").Value = "'6.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.57%  "

$ws.Range("D37").Value = "0.0₃0794"
$ws.Range("E37").Value = "  +15.73%  "

$ws.Range("D38").Value = "'55.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.73%  "

$ws.Range("D39").Value = "'492.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.00%  "

$ws.Range("D40").Value = "'3.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.32%  "

$ws.Range("D41").Value = "'0.0419"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.94%  "

$ws.Range("E42").Value = "  +2.35%  "

$ws.Range("D43").Value = "'8.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.29%  "

$ws.Range("D44").Value = "'0.291"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.37%  "

$ws.Range("D45").Value = "2.944.22"
$ws.Range("E45").Value = "  +2.58%  "

$ws.Range("D46").Value = "'2.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.13%  "

$ws.Range("D47").Value = "'28.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.11%  "

$ws.Range("E48").Value = "  -0.06%  "

$ws.Range("E49").Value = "  +1.08%  "

$ws.Range("D51").Value = "'2.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.72%  "
